$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set the new value in B5 (adds a new shared string "Redirecting to Plant")
$ws.Range("B5").Value = "Redirecting to Plant"

# Update the sheet's selected cell to B5 (matches <selection activeCell="B5" sqref="B5"/>)
$ws.Activate()
$ws.Range("B5").Select()
